$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.573649333333333
$ws.Range("H2").Value = 4.720948
$ws.Range("I2").Value = 0.162950296453897
$ws.Range("J2").Value = 0.1728167686459121
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.107177
$ws.Range("N2").Value = 0.321531
$ws.Range("O2").Value = 0.003526763356587491
$ws.Range("P2").Value = 0.003549676734010809
$ws.Range("Q2").Value = 0.1686590145986667
$ws.Range("R2").Value = 1.517931131388
$ws.Range("S2").Value = 0.0005746871344786726
$ws.Range("T2").Value = 0.0006134436629093226

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.573649333333333
$ws.Range("H3").Value = 4.720948
$ws.Range("I3").Value = 0.162950296453897
$ws.Range("J3").Value = 0.1728167686459121
$ws.Range("O3").Value = 0.9757678722356318
$ws.Range("P3").Value = 0.9821074349659524
$ws.Range("Q3").Value = 46.66376254048977
$ws.Range("R3").Value = 419.973862864408
$ws.Range("S3").Value = 0.1590016640509845
$ws.Range("T3").Value = 0.1697246333739411

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.573649333333333
$ws.Range("H4").Value = 4.720948
$ws.Range("I4").Value = 0.162950296453897
$ws.Range("J4").Value = 0.1728167686459121
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04072766666666667
$ws.Range("N4").Value = 0.122183
$ws.Range("O4").Value = 0.001340183457265176
$ws.Range("P4").Value = 0.001348890627627329
$ws.Range("Q4").Value = 0.06409106549822222
$ws.Range("R4").Value = 0.576819589484
$ws.Range("S4").Value = 0.0002183832916639691
$ws.Range("T4").Value = 0.0002331109195233112

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.573649333333333
$ws.Range("H5").Value = 4.720948
$ws.Range("I5").Value = 0.162950296453897
$ws.Range("J5").Value = 0.1728167686459121
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5885005
$ws.Range("N5").Value = 1.177001
$ws.Range("O5").Value = 0.01936518095051565
$ws.Range("P5").Value = 0.01299399767240936
$ws.Range("Q5").Value = 0.9260934194913333
$ws.Range("R5").Value = 5.556560516948
$ws.Range("S5").Value = 0.003155561976769885
$ws.Range("T5").Value = 0.002245580689538289

# Row 6
$ws.Range("I6").Value = 0.4548971409363525
$ws.Range("J6").Value = 0.4824406931050072
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.107177
$ws.Range("N6").Value = 0.321531
$ws.Range("O6").Value = 0.003526763356587491
$ws.Range("P6").Value = 0.003549676734010809
$ws.Range("Q6").Value = 0.4708337769473334
$ws.Range("R6").Value = 4.237503992526
$ws.Range("S6").Value = 0.001604314567670744
$ws.Range("T6").Value = 0.001712508503854893

# Row 7
$ws.Range("I7").Value = 0.4548971409363525
$ws.Range("J7").Value = 0.4824406931050072
$ws.Range("O7").Value = 0.9757678722356318
$ws.Range("P7").Value = 0.9821074349659524
$ws.Range("S7").Value = 0.443874015297537
$ws.Range("T7").Value = 0.4738085916285549

# Row 8
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.4548971409363525
$ws.Range("J8").Value = 0.4824406931050072
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04072766666666667
$ws.Range("N8").Value = 0.122183
$ws.Range("O8").Value = 0.001340183457265176
$ws.Range("P8").Value = 0.001348890627627329
$ws.Range("Q8").Value = 0.1789186217464445
$ws.Range("R8").Value = 1.610267595718
$ws.Range("S8").Value = 0.000609645623040125
$ws.Range("T8").Value = 0.0006507597293153766

# Row 9
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.4548971409363525
$ws.Range("J9").Value = 0.4824406931050072
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5885005
$ws.Range("N9").Value = 1.177001
$ws.Range("O9").Value = 0.01936518095051565
$ws.Range("P9").Value = 0.01299399767240936
$ws.Range("Q9").Value = 2.585311336857667
$ws.Range("R9").Value = 15.511868021146
$ws.Range("S9").Value = 0.00880916544810469
$ws.Range("T9").Value = 0.006268833243282024

# Row 10
$ws.Range("G10").Value = 1.149467
$ws.Range("H10").Value = 3.448401
$ws.Range("I10").Value = 0.11902651019285
$ws.Range("J10").Value = 0.1262334424813261
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.107177
$ws.Range("N10").Value = 0.321531
$ws.Range("O10").Value = 0.003526763356587491
$ws.Range("P10").Value = 0.003549676734010809
$ws.Range("Q10").Value = 0.123196424659
$ws.Range("R10").Value = 1.108767821931
$ws.Range("S10").Value = 0.0004197783346106309
$ws.Range("T10").Value = 0.000448087913830055

# Row 11
$ws.Range("G11").Value = 1.149467
$ws.Range("H11").Value = 3.448401
$ws.Range("I11").Value = 0.11902651019285
$ws.Range("J11").Value = 0.1262334424813261
$ws.Range("O11").Value = 0.9757678722356318
$ws.Range("P11").Value = 0.9821074349659524
$ws.Range("Q11").Value = 34.08539246956066
$ws.Range("R11").Value = 306.7685322260459
$ws.Range("S11").Value = 0.11614224459051
$ws.Range("T11").Value = 0.1239748024022573

# Row 12
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 1.149467
$ws.Range("H12").Value = 3.448401
$ws.Range("I12").Value = 0.11902651019285
$ws.Range("J12").Value = 0.1262334424813261
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04072766666666667
$ws.Range("N12").Value = 0.122183
$ws.Range("O12").Value = 0.001340183457265176
$ws.Range("P12").Value = 0.001348890627627329
$ws.Range("Q12").Value = 0.04681510882033333
$ws.Range("R12").Value = 0.4213359793829999
$ws.Range("S12").Value = 0.0001595173599364625
$ws.Range("T12").Value = 0.0001702751074561943

# Row 13
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 1.149467
$ws.Range("H13").Value = 3.448401
$ws.Range("I13").Value = 0.11902651019285
$ws.Range("J13").Value = 0.1262334424813261
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5885005
$ws.Range("N13").Value = 1.177001
$ws.Range("O13").Value = 0.01936518095051565
$ws.Range("P13").Value = 0.01299399767240936
$ws.Range("Q13").Value = 0.6764619042334998
$ws.Range("R13").Value = 4.058771425400999
$ws.Range("S13").Value = 0.002304969907792937
$ws.Range("T13").Value = 0.001640277057782573

# Row 14
$ws.Range("G14").Value = 1.6540555
$ws.Range("H14").Value = 3.308111
$ws.Range("I14").Value = 0.1712762992154535
$ws.Range("J14").Value = 0.1210979348516435
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.107177
$ws.Range("N14").Value = 0.321531
$ws.Range("O14").Value = 0.003526763356587491
$ws.Range("P14").Value = 0.003549676734010809
$ws.Range("Q14").Value = 0.1772767063235
$ws.Range("R14").Value = 1.063660237941
$ws.Range("S14").Value = 0.0006040509759249761
$ws.Range("T14").Value = 0.0004298585218796357

# Row 15
$ws.Range("G15").Value = 1.6540555
$ws.Range("H15").Value = 3.308111
$ws.Range("I15").Value = 0.1712762992154535
$ws.Range("J15").Value = 0.1210979348516435
$ws.Range("O15").Value = 0.9757678722356318
$ws.Range("P15").Value = 0.9821074349659524
$ws.Range("Q15").Value = 49.04806391478434
$ws.Range("R15").Value = 294.288383488706
$ws.Range("S15").Value = 0.1671259100498564
$ws.Range("T15").Value = 0.1189311821768216

# Row 16
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 1.6540555
$ws.Range("H16").Value = 3.308111
$ws.Range("I16").Value = 0.1712762992154535
$ws.Range("J16").Value = 0.1210979348516435
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04072766666666667
$ws.Range("N16").Value = 0.122183
$ws.Range("O16").Value = 0.001340183457265176
$ws.Range("P16").Value = 0.001348890627627329
$ws.Range("Q16").Value = 0.06736582105216668
$ws.Range("R16").Value = 0.404194926313
$ws.Range("S16").Value = 0.0002295416628301512
$ws.Range("T16").Value = 0.0001633478693464068

# Row 17
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 1.6540555
$ws.Range("H17").Value = 3.308111
$ws.Range("I17").Value = 0.1712762992154535
$ws.Range("J17").Value = 0.1210979348516435
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5885005
$ws.Range("N17").Value = 1.177001
$ws.Range("O17").Value = 0.01936518095051565
$ws.Range("P17").Value = 0.01299399767240936
$ws.Range("Q17").Value = 0.97341248877775
$ws.Range("R17").Value = 3.893649955111
$ws.Range("S17").Value = 0.003316796526841919
$ws.Range("T17").Value = 0.001573546283595837

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.8870146666666666
$ws.Range("H18").Value = 2.661044
$ws.Range("I18").Value = 0.09184975320144682
$ws.Range("J18").Value = 0.097411160916111
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.107177
$ws.Range("N18").Value = 0.321531
$ws.Range("O18").Value = 0.003526763356587491
$ws.Range("P18").Value = 0.003549676734010809
$ws.Range("Q18").Value = 0.09506757092933334
$ws.Range("R18").Value = 0.8556081383640001
$ws.Range("S18").Value = 0.0003239323439024672
$ws.Range("T18").Value = 0.0003457781315369022

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.8870146666666666
$ws.Range("H19").Value = 2.661044
$ws.Range("I19").Value = 0.09184975320144682
$ws.Range("J19").Value = 0.097411160916111
$ws.Range("O19").Value = 0.9757678722356318
$ws.Range("P19").Value = 0.9821074349659524
$ws.Range("Q19").Value = 26.30283691449155
$ws.Range("R19").Value = 236.725532230424
$ws.Range("S19").Value = 0.08962403824674367
$ws.Range("T19").Value = 0.09566822538437741

# Row 20
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.8870146666666666
$ws.Range("H20").Value = 2.661044
$ws.Range("I20").Value = 0.09184975320144682
$ws.Range("J20").Value = 0.097411160916111
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.04072766666666667
$ws.Range("N20").Value = 0.122183
$ws.Range("O20").Value = 0.001340183457265176
$ws.Range("P20").Value = 0.001348890627627329
$ws.Range("Q20").Value = 0.03612603767244445
$ws.Range("R20").Value = 0.325134339052
$ws.Range("S20").Value = 0.0001230955197944682
$ws.Range("T20").Value = 0.0001313970019860397

# Row 21
$ws.Range("D21").Value = "MuSCs"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.8870146666666666
$ws.Range("H21").Value = 2.661044
$ws.Range("I21").Value = 0.09184975320144682
$ws.Range("J21").Value = 0.097411160916111
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.5885005
$ws.Range("N21").Value = 1.177001
$ws.Range("O21").Value = 0.01936518095051565
$ws.Range("P21").Value = 0.01299399767240936
$ws.Range("Q21").Value = 0.5220085748406667
$ws.Range("R21").Value = 3.132051449044
$ws.Range("S21").Value = 0.001778687091006222
$ws.Range("T21").Value = 0.00126576039821064
